$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 10: zigbee_join_cnt variable (written after first power-on post flash = network config)
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "zigbee_join_cnt"
$ws.Range("C10").Value = "u8"
$ws.Range("D10").Value = "烧写后第一次上电为配网"
$ws.Range("E10").Value = "0x2F81"

# Match the bordered/left-aligned look used by the rest of the table (row 9 -> row 10)
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column D widened slightly to fit the new description text
$ws.Columns.Item(4).ColumnWidth = 22.2

# Selection moves to D11 (just past the newly-added last row)
$ws.Range("D11").Select()
